# Updated 2D training schedules, no break screen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (trials 1-5), columns A-J (J is shared string "train_dim2_1")
$data = @(
    @(1, 3, 7, 7, 5, 4, -2, 23, 5),
    @(2, 2, 7, 7, 6, 5, -1, 12, 5),
    @(3, 1, 8, 2, 3, 1, -5, 56, 5),
    @(4, 4, 9, 6, 5, 2, -4, 45, 5),
    @(5, 2, 5, 5, 2, 3, -3, 34, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("H10").Select()
